$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 2)
        $text = $cell.Text
        if ($text -eq "level") {
            $cell.Value = "category"
        } elseif ($text -eq "level_1") {
            $cell.Value = "category_1"
        } elseif ($text -eq "level_2") {
            $cell.Value = "category_2"
        } elseif ($text -eq "level_3") {
            $cell.Value = "category_3"
        } elseif ($text -eq "level_4") {
            $cell.Value = "category_4"
        } elseif ($text -eq "level_5") {
            $cell.Value = "category_5"
        }
    }
}
